$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A6: fix text to introduce "Manejoa los errores" wording (revert of previous correction)
$ws.Range("A6").Value = "Cumple con Clean Architecture (Plural,Metodos CRUD, Manejoa los errores)"

# Set E7 to "Ok"
$ws.Range("E7").Value = "Ok"

# Update selection to E8
$ws.Range("E8").Select()
